{"js": "// Replace the three-digit x one-digit multiplication prompts in the\n// worksheet table with a new set of problems, matching the author's\n// regenerated output (commit \"Update master to output generated at\n// 9a8706d\"). Each old prompt string is unique in the document, so a\n// plain search-and-replace on the exact \"NNN\u00d7N=\" text is unambiguous.\n\nconst replacements = [\n  [\"296\u00d74=\", \"627\u00d72=\"],\n  [\"235\u00d77=\", \"980\u00d79=\"],\n  [\"940\u00d79=\", \"698\u00d72=\"],\n  [\"205\u00d73=\", \"594\u00d78=\"],\n  [\"588\u00d74=\", \"810\u00d78=\"],\n  [\"292\u00d74=\", \"155\u00d76=\"],\n  [\"518\u00d79=\", \"712\u00d76=\"],\n  [\"318\u00d74=\", \"469\u00d77=\"],\n  [\"131\u00d75=\", \"856\u00d78=\"],\n  [\"746\u00d78=\", \"183\u00d75=\"],\n  [\"464\u00d72=\", \"595\u00d77=\"],\n  [\"585\u00d78=\", \"392\u00d78=\"],\n  [\"193\u00d75=\", \"481\u00d73=\"],\n  [\"795\u00d76=\", \"168\u00d74=\"],\n  [\"873\u00d74=\", \"612\u00d78=\"],\n  [\"386\u00d75=\", \"960\u00d79=\"],\n  [\"438\u00d79=\", \"289\u00d79=\"],\n  [\"828\u00d74=\", \"488\u00d76=\"],\n  [\"489\u00d77=\", \"590\u00d74=\"],\n  [\"759\u00d74=\", \"367\u00d79=\"],\n  [\"384\u00d79=\", \"714\u00d79=\"],\n  [\"153\u00d78=\", \"453\u00d77=\"],\n  [\"507\u00d77=\", \"683\u00d75=\"],\n  [\"505\u00d77=\", \"627\u00d74=\"],\n  [\"464\u00d75=\", \"555\u00d73=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text to replace: ${oldText}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the three-digit x one-digit multiplication prompts in the\n# worksheet table with a new set of problems, matching the author's\n# regenerated output (commit \"Update master to output generated at\n# 9a8706d\"). Each old prompt string is unique in the document, so a\n# plain Find/Replace on the exact \"NNN\u00d7N=\" text is unambiguous.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"296\u00d74=\", \"627\u00d72=\"),\n    @(\"235\u00d77=\", \"980\u00d79=\"),\n    @(\"940\u00d79=\", \"698\u00d72=\"),\n    @(\"205\u00d73=\", \"594\u00d78=\"),\n    @(\"588\u00d74=\", \"810\u00d78=\"),\n    @(\"292\u00d74=\", \"155\u00d76=\"),\n    @(\"518\u00d79=\", \"712\u00d76=\"),\n    @(\"318\u00d74=\", \"469\u00d77=\"),\n    @(\"131\u00d75=\", \"856\u00d78=\"),\n    @(\"746\u00d78=\", \"183\u00d75=\"),\n    @(\"464\u00d72=\", \"595\u00d77=\"),\n    @(\"585\u00d78=\", \"392\u00d78=\"),\n    @(\"193\u00d75=\", \"481\u00d73=\"),\n    @(\"795\u00d76=\", \"168\u00d74=\"),\n    @(\"873\u00d74=\", \"612\u00d78=\"),\n    @(\"386\u00d75=\", \"960\u00d79=\"),\n    @(\"438\u00d79=\", \"289\u00d79=\"),\n    @(\"828\u00d74=\", \"488\u00d76=\"),\n    @(\"489\u00d77=\", \"590\u00d74=\"),\n    @(\"759\u00d74=\", \"367\u00d79=\"),\n    @(\"384\u00d79=\", \"714\u00d79=\"),\n    @(\"153\u00d78=\", \"453\u00d77=\"),\n    @(\"507\u00d77=\", \"683\u00d75=\"),\n    @(\"505\u00d77=\", \"627\u00d74=\"),\n    @(\"464\u00d75=\", \"555\u00d73=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1            # wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n\nWrite-Output \"done\"\n"}
